# Updates cryptos list values (price/volume) and reorders Cosmos/Toncoin rows
# as published by the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.191.42'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '2.483.90'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('D5').Formula = "'320.85"
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').Formula = "'108.02"
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Formula = "'0.539"
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Formula = "'39.03"
$ws.Range('E10').Value = '  +3.86%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Formula = "'18.35"
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Formula = "'7.16"
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '2.872.33'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '2.486.00'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '47.108.48'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').Formula = "'13.34"
$ws.Range('E19').Value = '  +4.97%  '
$ws.Range('D20').Formula = "'6.60"
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').Formula = "'2.73"
$ws.Range('E22').Value = '  +14.32%  '
$ws.Range('D23').Formula = "'70.45"
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('D24').Formula = "'245.72"
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('D25').Formula = "'2.54"
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Formula = "'25.66"
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Formula = "'2.29"
$ws.Range('E28').Value = '  +3.76%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Formula = "'9.96"
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D31').Formula = "'34.45"
$ws.Range('E31').Value = '  -2.14%  '
$ws.Range('D32').Formula = "'49.71"
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').Formula = "'20.22"
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('E34').Value = '  -1.01%  '
$ws.Range('D35').Formula = "'0.0780"
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E37').Value = '  +1.96%  '
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('D40').Formula = "'22.86"
$ws.Range('E40').Value = '  +8.07%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D43').Formula = "'118.75"
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '1.993.75'
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').Formula = "'2.01"
$ws.Range('E47').Value = '  -3.97%  '
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('D49').Formula = "'9.11"
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D51').Formula = "'56.51"
$ws.Range('E51').Value = '  +2.82%  '
